# Fruta / hortaliza, semanal
# A new weekly price-report row for "Vega Modelo de Temuco - Papaya" is
# inserted at row 80, pushing the existing rows 80:134 down to 81:135.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 80 (shifts 80:134 -> 81:135,
# mirroring Excel's Rows.Insert default of xlShiftDown).
$ws.Rows.Item(80).Insert()

# Populate the newly inserted row 80 with the new weekly observation.
$ws.Cells.Item(80, 1).Value = 10
$ws.Cells.Item(80, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(80, 3).Value = "La Araucanía"
$ws.Cells.Item(80, 4).Value = "11/6/2023"
$ws.Cells.Item(80, 5).Value = 9
$ws.Cells.Item(80, 6).Value = "Fruta"
$ws.Cells.Item(80, 7).Value = 100108
$ws.Cells.Item(80, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(80, 9).Value = 100108004
$ws.Cells.Item(80, 10).Value = "Papaya"
$ws.Cells.Item(80, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(80, 12).Value = "Primera"
$ws.Cells.Item(80, 13).Value = 50
$ws.Cells.Item(80, 14).Value = 2600
$ws.Cells.Item(80, 15).Value = 2600
$ws.Cells.Item(80, 16).Value = 2600
$ws.Cells.Item(80, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(80, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(80, 19).Value = 2600
$ws.Cells.Item(80, 20).Value = 1
